$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix E11 formatting: it was missing the green "action" highlight that the
# other rows with a filled "Action recommandée" cell use (style reused from E3,
# so no new style/fill gets created). ---
$ws.Range("E3").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New column G: running reference numbers added alongside existing rows
# (part of moving the JS <script> includes to the end of the HTML, which
# required re-numbering the recommended-action cross references). ---
$ws.Range("G3").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 5
$ws.Range("G9").Value = 6
$ws.Range("G11").Value = 10
$ws.Range("G12").Value = 6
$ws.Range("G13").Value = 7
$ws.Range("G15").Value = 7
$ws.Range("G17").Value = 7
$ws.Range("G18").Value = 8
$ws.Range("G19").Value = 8
$ws.Range("G20").Value = 8
$ws.Range("G21").Value = 8
$ws.Range("G24").Value = 1
$ws.Range("G26").Value = 4
$ws.Range("G30").Value = 9

# --- Move the active selection from E1 down to E3 ---
$ws.Range("E3").Select()
